$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new row at 377 ("AR Birdehide") - shifts the previous rows
#    377-381 down to 378-382, inheriting formatting (date style) the same
#    way Excel's own Insert does.
# ---------------------------------------------------------------------------
$ws.Rows.Item(377).Insert()

$ws.Cells.Item(377,1).Value = "AR Birdehide"
$ws.Cells.Item(377,2).Value = -35.465475320000003
$ws.Cells.Item(377,3).Value = 137.46179395999999
$ws.Cells.Item(377,4).Value = 46005
$ws.Cells.Item(377,5).Value = 0
$ws.Cells.Item(377,6).Value = 0
$ws.Cells.Item(377,7).Value = 0
$ws.Cells.Item(377,8).Value = 807
$ws.Cells.Item(377,9).Value = 0
$ws.Cells.Item(377,10).Value = 0
$ws.Cells.Item(377,11).Value = 15
$ws.Cells.Item(377,12).Value = 0
$ws.Cells.Item(377,13).Value = 0
$ws.Cells.Item(377,14).Value = 15956
$ws.Cells.Item(377,15).Value = 16778
$ws.Cells.Item(377,16).Value = "LEWK5222"

# ---------------------------------------------------------------------------
# 2) Append 5 brand-new rows (383-387) after the (now shifted) last row 382.
#    Pre-copy the date number format from the row above so column D keeps
#    the same style index as the rest of the table instead of Excel's
#    default "General" format.
# ---------------------------------------------------------------------------
$ws.Range("D382").Copy()
$ws.Range("D383:D387").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 383 - Island Beach
$ws.Cells.Item(383,1).Value = "Island Beach"
$ws.Cells.Item(383,2).Value = -35.788755899999998
$ws.Cells.Item(383,3).Value = 137.78784780000001
$ws.Cells.Item(383,4).Value = 46010
$ws.Cells.Item(383,5).Value = 111
$ws.Cells.Item(383,6).Value = 0
$ws.Cells.Item(383,7).Value = 888
$ws.Cells.Item(383,8).Value = 8111
$ws.Cells.Item(383,9).Value = 0
$ws.Cells.Item(383,10).Value = 111
$ws.Cells.Item(383,11).Value = 111
$ws.Cells.Item(383,12).Value = 0
$ws.Cells.Item(383,13).Value = 333
$ws.Cells.Item(383,14).Value = 111
$ws.Cells.Item(383,15).Value = 9776
$ws.Cells.Item(383,16).Value = "LEWK5222"

# Row 384 - Garden Island (note: no value at all in column J, left blank)
$ws.Cells.Item(384,1).Value = "Garden Island"
$ws.Cells.Item(384,2).Value = -34.804169999999999
$ws.Cells.Item(384,3).Value = 138.53970000000001
$ws.Cells.Item(384,4).Value = 46011
$ws.Cells.Item(384,5).Value = 0
$ws.Cells.Item(384,6).Value = 0.0275
$ws.Cells.Item(384,7).Value = 0.0050000000000000001
$ws.Cells.Item(384,8).Value = 0.22
$ws.Cells.Item(384,9).Value = 0
$ws.Cells.Item(384,11).Value = 1.2829999999999999
$ws.Cells.Item(384,12).Value = 0
$ws.Cells.Item(384,13).Value = 0
$ws.Cells.Item(384,14).Value = 1.393
$ws.Cells.Item(384,15).Value = 2.9284999999999997
$ws.Cells.Item(384,16).Value = "HANC5232"

# Row 385 - Semaphore (note: no value at all in column J, left blank)
$ws.Cells.Item(385,1).Value = "Semaphore"
$ws.Cells.Item(385,2).Value = -34.940480000000001
$ws.Cells.Item(385,3).Value = 138.49844100000001
$ws.Cells.Item(385,4).Value = 46011
$ws.Cells.Item(385,5).Value = 0
$ws.Cells.Item(385,6).Value = 0.03
$ws.Cells.Item(385,7).Value = 0
$ws.Cells.Item(385,8).Value = 9
$ws.Cells.Item(385,9).Value = 0.6
$ws.Cells.Item(385,11).Value = 8.1999999999999993
$ws.Cells.Item(385,12).Value = 0.025000000000000001
$ws.Cells.Item(385,13).Value = 0
$ws.Cells.Item(385,14).Value = 0.2
$ws.Cells.Item(385,15).Value = 18.054999999999996
$ws.Cells.Item(385,16).Value = "HANC5232"

# Row 386 - Island Beach
$ws.Cells.Item(386,1).Value = "Island Beach"
$ws.Cells.Item(386,2).Value = -35.788755899999998
$ws.Cells.Item(386,3).Value = 137.78784780000001
$ws.Cells.Item(386,4).Value = 46011
$ws.Cells.Item(386,5).Value = 111
$ws.Cells.Item(386,6).Value = 1222
$ws.Cells.Item(386,7).Value = 0
$ws.Cells.Item(386,8).Value = 34444
$ws.Cells.Item(386,9).Value = 0
$ws.Cells.Item(386,10).Value = 2222
$ws.Cells.Item(386,11).Value = 333
$ws.Cells.Item(386,12).Value = 0
$ws.Cells.Item(386,13).Value = 2000
$ws.Cells.Item(386,14).Value = 223
$ws.Cells.Item(386,15).Value = 40555
$ws.Cells.Item(386,16).Value = "LEWK5222"

# Row 387 - Brownlow
$ws.Cells.Item(387,1).Value = "Brownlow"
$ws.Cells.Item(387,2).Value = -35.671777499999997
$ws.Cells.Item(387,3).Value = 137.6137956
$ws.Cells.Item(387,4).Value = 46013
$ws.Cells.Item(387,5).Value = 222
$ws.Cells.Item(387,6).Value = 1222
$ws.Cells.Item(387,7).Value = 0
$ws.Cells.Item(387,8).Value = 1778
$ws.Cells.Item(387,9).Value = 0
$ws.Cells.Item(387,10).Value = 222
$ws.Cells.Item(387,11).Value = 6111
$ws.Cells.Item(387,12).Value = 666
$ws.Cells.Item(387,13).Value = 6889
$ws.Cells.Item(387,14).Value = 0
$ws.Cells.Item(387,15).Value = 17110
$ws.Cells.Item(387,16).Value = "LEWK5222"

# ---------------------------------------------------------------------------
# 3) Restore the frozen-pane view / selection roughly where the author left
#    it: top pane frozen at row 1, scrolled down near the new rows, with
#    A389 selected on the bottom pane.
# ---------------------------------------------------------------------------
$ws.Range("D1").Select()
$ws.Range("A389").Select()

Write-Host "done"
